$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.625.53'
$ws.Range("E2").Value = '  -1.84%  '

$ws.Range("D3").Value = '1.791.29'
$ws.Range("E3").Value = '  -0.25%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("E5").Value = '  -0.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '305.93'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.68%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4952'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3847'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.79%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09264'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +16.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.088'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.61%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '40.46'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.30%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.06%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.238'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.38'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.66%  '

$ws.Range("D15").Value = '1.788.30'
$ws.Range("E15").Value = '  -0.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.117'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.43%  '

$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.71'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.18%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001102'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.10%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06529'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.98'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.98%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.892'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.97%  '

$ws.Range("D23").Value = '27.671.46'
$ws.Range("E23").Value = '  -1.81%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.89'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.223'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.89%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.80'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.54%  '

$ws.Range("D27").Value = '1.997.48'
$ws.Range("E27").Value = '  -0.09%  '

$ws.Range("E28").Value = '  -0.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.380'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.05'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.62%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1065'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.27%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.047'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.59%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.605'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.83%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.490'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.90%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06776'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.800'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.29%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02286'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.53%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2114'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.29%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.30'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.97%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.885'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.59%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6093'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.87%  '

$ws.Range("E42").Value = '  +0.06%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.139'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.14%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.88'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.05%  '

$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5827'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.23%  '

$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.658'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.99%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.265'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.69%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.61'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.04%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.913'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.164'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.63%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06686'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.88%  '
